$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Taul1")
$ws.Activate()

# Row 117: "what I did" text renamed ("Notification tests" -> "Notification testit ...")
$ws.Range("C117").Value = "Notification testit, redux fake store, store Provider wrapper"

# Row 118: "Loginbar tests" -> "Loginbar testit"
$ws.Range("C118").Value = "Loginbar testit"

# New row 119: SignUpForm testit, 4 hours, client (copy formatting from the row above
# so B119 keeps the same centered number style as the rest of column B)
$ws.Range("B118").Copy()
$ws.Range("B119").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("B119").Value = 4
$ws.Range("C119").Value = "SignUpForm testit"
$ws.Range("D119").Value = "client"

# New row 120: just the trailing "client" marker in column D
$ws.Range("D120").Value = "client"

# Recompute the hours total so it also covers the newly added row
$ws.Range("B123").Formula = "=SUM(B2:B119)"

# Update the view/selection state to match what was recorded after the edit
$ws.Range("C119").Select()
$excel.ActiveWindow.ScrollRow = 100
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Left = 28680
$excel.ActiveWindow.Top = -120
$excel.ActiveWindow.Width = 29040
$excel.ActiveWindow.Height = 15840
